# Auto-generated edit script
# Applies quantity (F) corrections, recalculated Value (G), company
# Sub Total (B), and the cascading grand Sub Total / Grand Total rows,
# matching the target diff exactly (183 cell updates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 2183.28
$ws.Range("B12").Value = 4145.56
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 1024.5
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 384.2
$ws.Range("B32").Value = 9284.6
$ws.Range("F36").Value = 3
$ws.Range("G36").Value = 93.69
$ws.Range("F37").Value = 41
$ws.Range("G37").Value = 2871.64
$ws.Range("F39").Value = 79
$ws.Range("G39").Value = 585.39
$ws.Range("F46").Value = 9
$ws.Range("G46").Value = 351.36
$ws.Range("F50").Value = 35
$ws.Range("G50").Value = 3273.9
$ws.Range("F52").Value = 36
$ws.Range("G52").Value = 1605.6
$ws.Range("F55").Value = 21
$ws.Range("G55").Value = 741.72
$ws.Range("F57").Value = 2
$ws.Range("G57").Value = 118.98
$ws.Range("B61").Value = 25981.81
$ws.Range("F118").Value = 19
$ws.Range("G118").Value = 2368.16
$ws.Range("F122").Value = 90
$ws.Range("G122").Value = 9221.4
$ws.Range("F131").Value = 73
$ws.Range("G131").Value = 3452.9
$ws.Range("B133").Value = 207770.4
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 811.9299999999999
$ws.Range("B140").Value = 811.9299999999999
$ws.Range("F144").Value = 35
$ws.Range("G144").Value = 4950.05
$ws.Range("F147").Value = 33
$ws.Range("G147").Value = 3436.95
$ws.Range("B151").Value = 8893.59
$ws.Range("F166").Value = 22
$ws.Range("G166").Value = 1088.56
$ws.Range("B176").Value = 15211.89
$ws.Range("F200").Value = 3
$ws.Range("G200").Value = 250.5
$ws.Range("F213").Value = 10
$ws.Range("G213").Value = 723.5
$ws.Range("B216").Value = 9816.82
$ws.Range("F252").Value = 29
$ws.Range("G252").Value = 4161.5
$ws.Range("B259").Value = 18873.64
$ws.Range("F316").Value = 15
$ws.Range("G316").Value = 4000.95
$ws.Range("F337").Value = 35
$ws.Range("G337").Value = 10597.65
$ws.Range("F340").Value = 1
$ws.Range("G340").Value = 102.24
$ws.Range("F354").Value = 40
$ws.Range("G354").Value = 4044
$ws.Range("F356").Value = 12
$ws.Range("G356").Value = 2383.08
$ws.Range("B380").Value = 258085.89
$ws.Range("F435").Value = 75
$ws.Range("G435").Value = 4119
$ws.Range("B436").Value = 58047
$ws.Range("D436").Value = 105.54
$ws.Range("E436").Value = 126.1
$ws.Range("F436").Value = 62
$ws.Range("G436").Value = 6543.48
$ws.Range("B437").Value = 47097
$ws.Range("D437").Value = 112.28
$ws.Range("E437").Value = 134.16
$ws.Range("F437").Value = 15
$ws.Range("G437").Value = 1684.2
$ws.Range("F441").Value = 24
$ws.Range("G441").Value = 6036
$ws.Range("F444").Value = 24
$ws.Range("G444").Value = 590.16
$ws.Range("F446").Value = 40
$ws.Range("G446").Value = 6658.8
$ws.Range("B447").Value = 38470.24
$ws.Range("F450").Value = 54
$ws.Range("G450").Value = 2694.06
$ws.Range("F451").Value = 250
$ws.Range("G451").Value = 12955
$ws.Range("F462").Value = 41
$ws.Range("G462").Value = 7644.04
$ws.Range("F467").Value = 0
$ws.Range("G467").Value = 0
$ws.Range("B473").Value = 136849.99
$ws.Range("F491").Value = 476
$ws.Range("G491").Value = 6402.2
$ws.Range("F493").Value = 550
$ws.Range("G493").Value = 7045.5
$ws.Range("F496").Value = 306
$ws.Range("G496").Value = 5027.58
$ws.Range("F501").Value = 93
$ws.Range("G501").Value = 1809.78
$ws.Range("F504").Value = 950
$ws.Range("G504").Value = 6165.5
$ws.Range("F505").Value = 409
$ws.Range("G505").Value = 5378.35
$ws.Range("F506").Value = 325
$ws.Range("G506").Value = 8547.5
$ws.Range("F507").Value = 256
$ws.Range("G507").Value = 4206.08
$ws.Range("F508").Value = 550
$ws.Range("G508").Value = 8101.5
$ws.Range("B509").Value = 94912.96000000001
$ws.Range("F512").Value = 28
$ws.Range("G512").Value = 812.5599999999999
$ws.Range("B516").Value = 6638.12
$ws.Range("F558").Value = 602
$ws.Range("G558").Value = 11949.7
$ws.Range("B563").Value = 36937.82
$ws.Range("F574").Value = 12
$ws.Range("G574").Value = 312.96
$ws.Range("F575").Value = 17
$ws.Range("G575").Value = 887.0599999999999
$ws.Range("B584").Value = 23371.92
$ws.Range("F595").Value = 31
$ws.Range("G595").Value = 3053.5
$ws.Range("B601").Value = 64304.69
$ws.Range("F619").Value = 26
$ws.Range("G619").Value = 2700.62
$ws.Range("F622").Value = 24
$ws.Range("G622").Value = 1205.52
$ws.Range("F636").Value = 48
$ws.Range("G636").Value = 5895.84
$ws.Range("B640").Value = 208550.55
$ws.Range("F669").Value = 115
$ws.Range("G669").Value = 1823.9
$ws.Range("F670").Value = 73
$ws.Range("G670").Value = 2417.03
$ws.Range("B677").Value = 20981.96
$ws.Range("F679").Value = 27
$ws.Range("G679").Value = 6062.31
$ws.Range("F681").Value = 24
$ws.Range("G681").Value = 1917.12
$ws.Range("F682").Value = 20
$ws.Range("G682").Value = 1640.8
$ws.Range("F692").Value = 27
$ws.Range("G692").Value = 2390.04
$ws.Range("F693").Value = 14
$ws.Range("G693").Value = 1390.34
$ws.Range("F694").Value = 19
$ws.Range("G694").Value = 2092.09
$ws.Range("B695").Value = 47857.65
$ws.Range("F708").Value = 5
$ws.Range("G708").Value = 4733.55
$ws.Range("B716").Value = 108516.44
$ws.Range("F755").Value = 238
$ws.Range("G755").Value = 19411.28
$ws.Range("F758").Value = 261
$ws.Range("G758").Value = 34060.5
$ws.Range("F761").Value = 34
$ws.Range("G761").Value = 3792.36
$ws.Range("F763").Value = 117
$ws.Range("G763").Value = 2541.24
$ws.Range("F765").Value = 75
$ws.Range("G765").Value = 6801
$ws.Range("F771").Value = 503
$ws.Range("G771").Value = 67910.03
$ws.Range("F773").Value = 579
$ws.Range("G773").Value = 69891.09
$ws.Range("B775").Value = 254893.77
$ws.Range("F791").Value = 67
$ws.Range("G791").Value = 5400.2
$ws.Range("B793").Value = 13433.45
$ws.Range("F800").Value = 12
$ws.Range("G800").Value = 448.8
$ws.Range("B801").Value = 642.66
$ws.Range("F852").Value = 697
$ws.Range("G852").Value = 21070.31
$ws.Range("F853").Value = 3174
$ws.Range("G853").Value = 517711.14
$ws.Range("F855").Value = 232
$ws.Range("G855").Value = 33558.8
$ws.Range("F856").Value = 125
$ws.Range("G856").Value = 4767.5
$ws.Range("B861").Value = 636801.23
$ws.Range("B867").Value = 3536732.29
$ws.Range("B868").Value = 3536732.29
